$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.01"
$ws.Range("E2").Value = "'6.47%"
$ws.Range("D3").Value = "'31.84"
$ws.Range("E3").Value = "'8.37%"
$ws.Range("D4").Value = "'5.256"
$ws.Range("E4").Value = "'3.78%"
$ws.Range("D5").Value = "'0.07508"
$ws.Range("E5").Value = "'11.91%"
$ws.Range("D6").Value = "'7.825"
$ws.Range("E6").Value = "'6.98%"
$ws.Range("D7").Value = "'3.761"
$ws.Range("E7").Value = "'9.35%"
$ws.Range("D8").Value = "'1.482"
$ws.Range("E8").Value = "'7.26%"
$ws.Range("D9").Value = "'0.9147"
$ws.Range("E9").Value = "'1.62%"
$ws.Range("D10").Value = "'0.01654"
$ws.Range("E10").Value = "'2,453.74%"
$ws.Range("D11").Value = "'0.1691"
$ws.Range("E11").Value = "'7.81%"
$ws.Range("D12").Value = "'0.07540"
$ws.Range("E12").Value = "'6.13%"
$ws.Range("E13").Value = "'5.79%"
$ws.Range("D14").Value = "'0.02996"
$ws.Range("E14").Value = "'2.65%"
$ws.Range("D15").Value = "'0.09913"
$ws.Range("E15").Value = "'10.27%"
$ws.Range("D16").Value = "'0.001491"
$ws.Range("E16").Value = "'-5.52%"
$ws.Range("D17").Value = "'0.04557"
$ws.Range("E17").Value = "'1.57%"
$ws.Range("D18").Value = "'0.006361"
$ws.Range("E18").Value = "'2.28%"
$ws.Range("D19").Value = "'3.496"
$ws.Range("E19").Value = "'1.36%"
$ws.Range("D20").Value = "'2.231"
$ws.Range("E20").Value = "'0.16%"
$ws.Range("D21").Value = "'0.3311"
$ws.Range("E21").Value = "'2.43%"
$ws.Range("E22").Value = "'1.93%"
$ws.Range("D23").Value = "'4.473"
$ws.Range("E23").Value = "'13.68%"
$ws.Range("D25").Value = "'0.001214"
$ws.Range("E25").Value = "'1.04%"
$ws.Range("D26").Value = "'0.004442"
$ws.Range("E26").Value = "'1.69%"
$ws.Range("D27").Value = "'0.0001397"
$ws.Range("E27").Value = "'19.36%"
$ws.Range("D28").Value = "'0.0001740"
$ws.Range("E28").Value = "'7.52%"
$ws.Range("D40").Value = "'0.04498"
$ws.Range("E40").Value = "'6.16%"
$ws.Range("D41").Value = "'0.007217"
$ws.Range("E41").Value = "'6.45%"
$ws.Range("E42").Value = "'8.83%"
$ws.Range("D43").Value = "'0.002245"
$ws.Range("E43").Value = "'0.65%"
$ws.Range("D44").Value = "'0.01292"
$ws.Range("E44").Value = "'1.85%"
$ws.Range("D45").Value = "'0.00006220"
$ws.Range("E45").Value = "'8.08%"
$ws.Range("D46").Value = "'0.7091"
$ws.Range("E46").Value = "'-63.88%"
$ws.Range("D47").Value = "'0.01299"
$ws.Range("E47").Value = "'-13.44%"
